$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 43
$ws.Range("H43").Value = 50000500
$ws.Range("I43").Value = 50000500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 50000500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -50000431
$ws.Range("N43").ClearContents()

# ALC row 96
$ws.Range("H96").Value = 780
$ws.Range("I96").Value = 137.75
$ws.Range("K96").Value = 413.25
$ws.Range("M96").Value = 959.75

# ALC row 98
$ws.Range("H98").Value = 532.25
$ws.Range("I98").Value = 506.85715
$ws.Range("K98").Value = 506.85715
$ws.Range("M98").Value = 991.14285

# ALC row 122
$ws.Range("H122").Value = 532.25
$ws.Range("I122").Value = 506.85715
$ws.Range("K122").Value = 1520.57145
$ws.Range("M122").Value = 929.4285500000001

# ALC row 132
$ws.Range("H132").Value = 12343.368
$ws.Range("I132").Value = 12974.944
$ws.Range("K132").Value = 38924.83199999999
$ws.Range("M132").Value = -36394.83199999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 1165.5238
$ws.Range("I32").Value = 560.6286
$ws.Range("K32").Value = 560.6286
$ws.Range("M32").Value = -273.6286

# ARM row 45
$ws.Range("H45").Value = 3777.8
$ws.Range("J45").Value = 4296.3335
$ws.Range("L45").Value = 4296.3335
$ws.Range("N45").Value = -5050.3335

# ARM row 74
$ws.Range("H74").Value = 4649.1
$ws.Range("I74").Value = 4649.1
$ws.Range("K74").Value = 4649.1
$ws.Range("M74").Value = -3775.1

# ARM row 77
$ws.Range("H77").Value = 4649.1
$ws.Range("I77").Value = 4649.1
$ws.Range("K77").Value = 23245.5
$ws.Range("M77").Value = -18877.5

# ARM row 122
$ws.Range("H122").Value = 3372.25
$ws.Range("I122").Value = 2006
$ws.Range("K122").Value = 6018
$ws.Range("M122").Value = -3568

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 2199.8
$ws.Range("I20").Value = 1874.75
$ws.Range("K20").Value = 1874.75
$ws.Range("M20").Value = -1627.75

$ws = $wb.Worksheets.Item("CRP")
# CRP row 70
$ws.Range("H70").Value = 19999
$ws.Range("I70").Value = 19999
$ws.Range("K70").Value = 19999
$ws.Range("M70").Value = -19684

# CRP row 73
$ws.Range("H73").Value = 19999
$ws.Range("I73").Value = 19999
$ws.Range("K73").Value = 19999
$ws.Range("M73").Value = -18907

# CRP row 132
$ws.Range("H132").Value = 3103.8823
$ws.Range("I132").Value = 1420.4
$ws.Range("J132").Value = 5508.857
$ws.Range("K132").Value = 4261.200000000001
$ws.Range("L132").Value = 16526.571
$ws.Range("M132").Value = -1731.200000000001
$ws.Range("N132").Value = -21586.571

$ws = $wb.Worksheets.Item("CUL")
# CUL row 86
$ws.Range("H86").Value = 499.5
$ws.Range("J86").Value = 499.5
$ws.Range("L86").Value = 1498.5
$ws.Range("N86").Value = -3870.5

# CUL row 89
$ws.Range("H89").Value = 499.5
$ws.Range("J89").Value = 499.5
$ws.Range("L89").Value = 4495.5
$ws.Range("N89").Value = -16351.5

# CUL row 103
$ws.Range("H103").Value = 659.6667
$ws.Range("I103").Value = 292.83334
$ws.Range("K103").Value = 878.5000200000001
$ws.Range("M103").Value = 0.4999799999999368

# CUL row 107
$ws.Range("H107").Value = 367.66666
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# CUL row 129
$ws.Range("H129").Value = 1444
$ws.Range("I129").Value = 652
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 1956
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 3044
$ws.Range("N129").Value = -17500

# CUL row 131
$ws.Range("H131").Value = 1687.5
$ws.Range("I131").Value = 853.75
$ws.Range("J131").Value = 2243.3333
$ws.Range("K131").Value = 2561.25
$ws.Range("L131").Value = 6729.999899999999
$ws.Range("M131").Value = 2478.75
$ws.Range("N131").Value = -16809.9999

# CUL row 133
$ws.Range("H133").Value = 2000
$ws.Range("I133").Value = 2000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 6000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -940
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# GSM row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# GSM row 80
$ws.Range("H80").Value = 3181.3333
$ws.Range("J80").Value = 3844
$ws.Range("L80").Value = 3844
$ws.Range("N80").Value = -5840

# GSM row 83
$ws.Range("H83").Value = 3181.3333
$ws.Range("J83").Value = 3844
$ws.Range("L83").Value = 19220
$ws.Range("N83").Value = -29204

# GSM row 97
$ws.Range("H97").Value = 798.8333
$ws.Range("I97").Value = 698.3333
$ws.Range("J97").Value = 899.3333
$ws.Range("K97").Value = 698.3333
$ws.Range("L97").Value = 899.3333
$ws.Range("M97").Value = -202.3333
$ws.Range("N97").Value = -1891.3333

# GSM row 132
$ws.Range("H132").Value = 6080.7144
$ws.Range("I132").Value = 5513.2
$ws.Range("K132").Value = 16539.6
$ws.Range("M132").Value = -14009.6

$ws = $wb.Worksheets.Item("LTW")
# LTW row 20
$ws.Range("H20").Value = 128499.75
$ws.Range("J20").Value = 512499
$ws.Range("L20").Value = 512499
$ws.Range("N20").Value = -512951

# LTW row 22
$ws.Range("H22").Value = 3309.8
$ws.Range("I22").Value = 3162.5
$ws.Range("K22").Value = 3162.5
$ws.Range("M22").Value = -2867.5

# LTW row 27
$ws.Range("H27").Value = 3309.8
$ws.Range("I27").Value = 3162.5
$ws.Range("K27").Value = 3162.5
$ws.Range("M27").Value = -3055.5

# LTW row 42
$ws.Range("H42").Value = 28999.666
$ws.Range("J42").Value = 39999
$ws.Range("L42").Value = 39999
$ws.Range("N42").Value = -41125

# LTW row 46
$ws.Range("H46").Value = 6225
$ws.Range("I46").Value = 4962.5
$ws.Range("J46").Value = 8750
$ws.Range("K46").Value = 4962.5
$ws.Range("L46").Value = 8750
$ws.Range("M46").Value = -4774.5
$ws.Range("N46").Value = -9126

# LTW row 49
$ws.Range("H49").Value = 28999.666
$ws.Range("J49").Value = 39999
$ws.Range("L49").Value = 39999
$ws.Range("N49").Value = -40293

# LTW row 55
$ws.Range("H55").Value = 1393
$ws.Range("I55").Value = 1892.25
$ws.Range("J55").Value = 394.5
$ws.Range("K55").Value = 1892.25
$ws.Range("L55").Value = 394.5
$ws.Range("M55").Value = -1719.25
$ws.Range("N55").Value = -740.5

# LTW row 82
$ws.Range("H82").Value = 3050.125
$ws.Range("I82").Value = 851
$ws.Range("K82").Value = 851
$ws.Range("M82").Value = -490

# LTW row 85
$ws.Range("H85").Value = 3050.125
$ws.Range("I85").Value = 851
$ws.Range("K85").Value = 851
$ws.Range("M85").Value = 397

# LTW row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132
$ws.Range("H132").Value = 2780.4167
$ws.Range("I132").Value = 2986.6
$ws.Range("J132").Value = 1749.5
$ws.Range("K132").Value = 8959.799999999999
$ws.Range("L132").Value = 5248.5
$ws.Range("M132").Value = -6429.799999999999
$ws.Range("N132").Value = -10308.5
